$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "44.581.26", "0.0763") in the
# source workbook. Excel's COM layer auto-detects numeric-looking strings and
# would silently coerce them to real numbers (losing the text type + exact
# formatting, e.g. "0.0763" -> 7.63E-2). Force the whole Price column to a
# text number-format before writing so every assignment below is stored
# verbatim as a string, then restore the default "Normal" style afterwards
# so we do not leave a stray text-format style on cells that did not need it.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "44.581.26"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "2.434.74"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "311.66"
$ws.Range("E5").Value = "  +3.45%  "
$ws.Range("D6").Value = "101.91"
$ws.Range("E6").Value = "  +6.28%  "
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").Value = "35.51"
$ws.Range("E10").Value = "  +4.17%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "18.76"
$ws.Range("E13").Value = "  +2.70%  "
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "2.812.84"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "2.425.07"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").Value = "44.531.87"
$ws.Range("E18").Value = "  +3.93%  "
$ws.Range("D19").Value = "12.42"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").Value = "6.43"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "68.88"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "25.28"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("E28").Value = "  -4.14%  "
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").Value = "33.52"
$ws.Range("E30").Value = "  +5.98%  "
$ws.Range("D31").Value = "48.72"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("E32").Value = "  +15.83%  "
$ws.Range("D33").Value = "19.51"
$ws.Range("E33").Value = "  +12.40%  "
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").Value = "0.0763"
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("D38").Value = "4.55"
$ws.Range("E38").Value = "  +4.49%  "
$ws.Range("D39").Value = "2.92"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("D40").Value = "124.68"
$ws.Range("E40").Value = "  +6.27%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "21.93"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  -5.10%  "
$ws.Range("D44").Value = "0.0290"
$ws.Range("E44").Value = "  +3.67%  "
$ws.Range("D45").Value = "1.948.46"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").Value = "2.98"
$ws.Range("E47").Value = "  +8.93%  "
$ws.Range("D48").Value = "9.84"
$ws.Range("E48").Value = "  +7.13%  "
$ws.Range("D49").Value = "1.68"
$ws.Range("E49").Value = "  +10.73%  "
$ws.Range("D50").Value = "53.48"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("D51").Value = "73.90"
$ws.Range("E51").Value = "  +2.70%  "


# Restore default styling on the Price column (values already committed as
# text above; this just clears the temporary "@" text format we applied).
$priceCol.Style = "Normal"
